# This script applies the renaming/title-casing edits described in the commit,
# plus trims the trailing footer/metadata rows and fixes the sheet dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns (A1:D1) to the new short machine-readable names
# 2) Title-case the Spanish linking particles (de/del/la/las/los/el/y) inside
#    state/municipality names, and fix 'MonteMorelos' -> 'Montemorelos'
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B6').Value = 'Pabellón De Arteaga'
$ws.Range('B7').Value = 'Rincón De Romos'
$ws.Range('B8').Value = 'San Francisco De Los Romo'
$ws.Range('B29').Value = 'Amatenango De La Frontera'
$ws.Range('B32').Value = 'Bejucal De Ocampo'
$ws.Range('B34').Value = 'Benemérito De Las Américas'
$ws.Range('B39').Value = 'Comitán De Domínguez'
$ws.Range('B58').Value = 'Marqués De Comillas'
$ws.Range('B59').Value = 'Mazapa De Madero'
$ws.Range('B62').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B69').Value = 'Salto De Agua'
$ws.Range('B70').Value = 'San Cristóbal De Las Casas'
$ws.Range('B96').Value = 'Hidalgo Del Parral'
$ws.Range('B120').Value = 'San Juan De Sabinas'
$ws.Range('A128').Value = 'Ciudad De México'
$ws.Range('B154').Value = 'Pánuco De Coronado'
$ws.Range('B158').Value = 'San Pedro Del Gallo'
$ws.Range('A166').Value = 'Estado De México'
$ws.Range('B166').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B168').Value = 'Almoloya De Juárez'
$ws.Range('B169').Value = 'Almoloya Del Río'
$ws.Range('B172').Value = 'Atizapán De Zaragoza'
$ws.Range('B177').Value = 'Chapa De Mota'
$ws.Range('B179').Value = 'Coacalco De Berriozábal'
$ws.Range('B183').Value = 'Ecatepec De Morelos'
$ws.Range('B188').Value = 'Ixtapan De La Sal'
$ws.Range('B189').Value = 'Ixtapan Del Oro'
$ws.Range('B198').Value = 'Naucalpan De Juárez'
$ws.Range('B203').Value = 'San Felipe Del Progreso'
$ws.Range('B217').Value = 'Tlalnepantla De Baz'
$ws.Range('B223').Value = 'Valle De Bravo'
$ws.Range('B224').Value = 'Villa De Allende'
$ws.Range('B225').Value = 'Villa Del Carbón'
$ws.Range('B236').Value = 'Apaseo El Alto'
$ws.Range('B237').Value = 'Apaseo El Grande'
$ws.Range('B243').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B247').Value = 'Jaral Del Progreso'
$ws.Range('B256').Value = 'San Diego De La Unión'
$ws.Range('B258').Value = 'San Francisco Del Rincón'
$ws.Range('B260').Value = 'San Luis De La Paz'
$ws.Range('B261').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B262').Value = 'Silao De La Victoria'
$ws.Range('B264').Value = 'Valle De Santiago'
$ws.Range('B270').Value = 'Acapulco De Juárez'
$ws.Range('B273').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B274').Value = 'Alcozauca De Guerrero'
$ws.Range('B277').Value = 'Atlamajalcingo Del Monte'
$ws.Range('B278').Value = 'Atoyac De Álvarez'
$ws.Range('B279').Value = 'Ayutla De Los Libres'
$ws.Range('B282').Value = 'Buenavista De Cuéllar'
$ws.Range('B283').Value = 'Chilapa De Álvarez'
$ws.Range('B284').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B285').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B289').Value = 'Coyuca De Benítez'
$ws.Range('B290').Value = 'Coyuca De Catalán'
$ws.Range('B293').Value = 'Cuetzala Del Progreso'
$ws.Range('B294').Value = 'Cutzamala De Pinzón'
$ws.Range('B300').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B301').Value = 'Iguala De La Independencia'
$ws.Range('B303').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B304').Value = 'Zihuatanejo De Azueta'
$ws.Range('B306').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B320').Value = 'Taxco De Alarcón'
$ws.Range('B322').Value = 'Técpan De Galeana'
$ws.Range('B324').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B325').Value = 'Tixtla De Guerrero'
$ws.Range('B329').Value = 'Tlapa De Comonfort'
$ws.Range('B342').Value = 'Atotonilco De Tula'
$ws.Range('B343').Value = 'Atotonilco El Grande'
$ws.Range('B351').Value = 'Huejutla De Reyes'
$ws.Range('B354').Value = 'Jacala De Ledezma'
$ws.Range('B357').Value = 'Mineral Del Monte'
$ws.Range('B358').Value = 'Mixquiahuala De Juárez'
$ws.Range('B359').Value = 'Molango De Escamilla'
$ws.Range('B361').Value = 'Nopala De Villagrán'
$ws.Range('B362').Value = 'Omitlán De Juárez'
$ws.Range('B363').Value = 'Pachuca De Soto'
$ws.Range('B366').Value = 'Progreso De Obregón'
$ws.Range('B370').Value = 'Santiago De Anaya'
$ws.Range('B374').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B375').Value = 'Tezontepec De Aldama'
$ws.Range('B378').Value = 'Tula De Allende'
$ws.Range('B379').Value = 'Tulancingo De Bravo'
$ws.Range('B381').Value = 'Zacualtipán De Ángeles'
$ws.Range('B382').Value = 'Zapotlán De Juárez'
$ws.Range('B386').Value = 'Acatlán De Juárez'
$ws.Range('B387').Value = 'Ahualulco De Mercado'
$ws.Range('B398').Value = 'Encarnación De Díaz'
$ws.Range('B403').Value = 'Jilotlán De Los Dolores'
$ws.Range('B406').Value = 'Lagos De Moreno'
$ws.Range('B413').Value = 'San Diego De Alejandría'
$ws.Range('B414').Value = 'San Juan De Los Lagos'
$ws.Range('B416').Value = 'San Martín De Bolaños'
$ws.Range('B418').Value = 'San Miguel El Alto'
$ws.Range('B419').Value = 'Santa María Del Oro'
$ws.Range('B420').Value = 'Talpa De Allende'
$ws.Range('B421').Value = 'Tamazula De Gordiano'
$ws.Range('B422').Value = 'Teocuitatlán De Corona'
$ws.Range('B423').Value = 'Tepatitlán De Morelos'
$ws.Range('B425').Value = 'Tizapán El Alto'
$ws.Range('B431').Value = 'Unión De Tula'
$ws.Range('B432').Value = 'Valle De Guadalupe'
$ws.Range('B434').Value = 'Yahualica De González Gallo'
$ws.Range('B435').Value = 'Zacoalco De Torres'
$ws.Range('B437').Value = 'Zapotlán Del Rey'
$ws.Range('B438').Value = 'Zapotlán El Grande'
$ws.Range('B456').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B508').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B528').Value = 'Coatlán Del Río'
$ws.Range('B536').Value = 'Puente De Ixtla'
$ws.Range('B539').Value = 'Tlaltizapán De Zapata'
$ws.Range('B551').Value = 'Santa María Del Oro'
$ws.Range('B565').Value = 'Mier Y Noriega'
$ws.Range('B566').Value = 'Montemorelos'
$ws.Range('B569').Value = 'San Nicolás De Los Garza'
$ws.Range('B572').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B574').Value = 'Ayoquezco De Aldama'
$ws.Range('B577').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B578').Value = 'Chiquihuitlán De Benito Juárez'
$ws.Range('B582').Value = 'El Barrio De La Soledad'
$ws.Range('B583').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B584').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B585').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B586').Value = 'Huautla De Jiménez'
$ws.Range('B587').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B591').Value = 'Mariscala De Juárez'
$ws.Range('B592').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B594').Value = 'Oaxaca De Juárez'
$ws.Range('B595').Value = 'Ocotlán De Morelos'
$ws.Range('B597').Value = 'Putla Villa De Guerrero'
$ws.Range('B598').Value = 'Reforma De Pineda'
$ws.Range('B608').Value = 'San Francisco Del Mar'
$ws.Range('B619').Value = 'San Juan De Los Cués'
$ws.Range('B620').Value = 'San Juan Del Estado'
$ws.Range('B621').Value = 'San Juan Del Río'
$ws.Range('B639').Value = 'San Miguel El Grande'
$ws.Range('B692').Value = 'Santo Domingo De Morelos'
$ws.Range('B698').Value = 'Tataltepec De Valdés'
$ws.Range('B699').Value = 'Teotitlán Del Valle'
$ws.Range('B700').Value = 'Tepelmeme Villa De Morelos'
$ws.Range('B701').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B702').Value = 'Tlacolula De Matamoros'
$ws.Range('B703').Value = 'Totontepec Villa De Morelos'
$ws.Range('B705').Value = 'Villa De Etla'
$ws.Range('B706').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B707').Value = 'Villa Sola De Vega'
$ws.Range('B708').Value = 'Yutanduchi De Guerrero'
$ws.Range('B721').Value = 'Ayotoxco De Guerrero'
$ws.Range('B722').Value = 'Chalchicomula De Sesma'
$ws.Range('B731').Value = 'Cuetzalan Del Progreso'
$ws.Range('B740').Value = 'Izúcar De Matamoros'
$ws.Range('B744').Value = 'Los Reyes De Juárez'
$ws.Range('B746').Value = 'Palmar De Bravo'
$ws.Range('B756').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B758').Value = 'San Salvador El Seco'
$ws.Range('B759').Value = 'San Salvador El Verde'
$ws.Range('B764').Value = 'Tepango De Rodríguez'
$ws.Range('B765').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B767').Value = 'Tetela De Ocampo'
$ws.Range('B770').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B785').Value = 'Amealco De Bonfil'
$ws.Range('B787').Value = 'Cadereyta De Montes'
$ws.Range('B791').Value = 'Jalpan De Serra'
$ws.Range('B793').Value = 'Pinal De Amoles'
$ws.Range('B796').Value = 'San Juan Del Río'
$ws.Range('B806').Value = 'Ciudad Del Maíz'
$ws.Range('B814').Value = 'Mexquitic De Carmona'
$ws.Range('B821').Value = 'Santa María Del Río'
$ws.Range('B828').Value = 'Tanquián De Escobedo'
$ws.Range('B830').Value = 'Villa De Arista'
$ws.Range('B831').Value = 'Villa De Arriaga'
$ws.Range('B832').Value = 'Villa De La Paz'
$ws.Range('B833').Value = 'Villa De Ramos'
$ws.Range('B834').Value = 'Villa De Reyes'
$ws.Range('B867').Value = 'Jalpa De Méndez'
$ws.Range('B893').Value = 'Soto La Marina'
$ws.Range('B900').Value = 'Apetatitlán De Antonio Carvajal'
$ws.Range('B906').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B907').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B924').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B926').Value = 'Amatlán De Los Reyes'
$ws.Range('B935').Value = 'Camarón De Tejeda'
$ws.Range('B938').Value = 'Castillo De Teayo'
$ws.Range('B940').Value = 'Cazones De Herrera'
$ws.Range('B948').Value = 'Cosamaloapan De Carpio'
$ws.Range('B965').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B966').Value = 'Ixhuatlán De Madero'
$ws.Range('B967').Value = 'Ixhuatlán Del Café'
$ws.Range('B968').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B976').Value = 'Juchique De Ferrer'
$ws.Range('B980').Value = 'Lerdo De Tejada'
$ws.Range('B984').Value = 'Martínez De La Torre'
$ws.Range('B985').Value = 'Medellín De Bravo'
$ws.Range('B988').Value = 'Mixtla De Altamirano'
$ws.Range('B994').Value = 'Ozuluama De Mascareñas'
$ws.Range('B997').Value = 'Paso De Ovejas'
$ws.Range('B998').Value = 'Paso Del Macho'
$ws.Range('B1001').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1007').Value = 'Sayula De Alemán'
$ws.Range('B1009').Value = 'Soledad De Doblado'
$ws.Range('B1021').Value = 'Tlacotepec De Mejía'
$ws.Range('B1033').Value = 'Vega De Alatorre'
$ws.Range('B1061').Value = 'Nochistlán De Mejía'
$ws.Range('B1070').Value = 'Villa De Cos'

# 3) Remove the trailing sample-size / source / author / date footer rows
#    (rows 1076-1080); this also shrinks the used range / dimension accordingly
$ws.Range("A1076:A1080").EntireRow.Delete()

Write-Host "Applied header rename, particle title-casing, and footer row cleanup."
